# Update "想去人数" (interest counts) in column F for the "展览" and
# "全部类型" worksheets, matching a refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 627
$ws1.Range("F4").Value  = 1302
$ws1.Range("F5").Value  = 1170
$ws1.Range("F6").Value  = 14362
$ws1.Range("F7").Value  = 16658
$ws1.Range("F8").Value  = 15
$ws1.Range("F9").Value  = 107
$ws1.Range("F10").Value = 19
$ws1.Range("F11").Value = 50
$ws1.Range("F13").Value = 26
$ws1.Range("F14").Value = 51
$ws1.Range("F19").Value = 108
$ws1.Range("F20").Value = 38
$ws1.Range("F21").Value = 1271
$ws1.Range("F22").Value = 137
$ws1.Range("F24").Value = 42
$ws1.Range("F25").Value = 22
$ws1.Range("F27").Value = 6781
$ws1.Range("F29").Value = 23
$ws1.Range("F30").Value = 1125
$ws1.Range("F33").Value = 5766
$ws1.Range("F34").Value = 112
$ws1.Range("F35").Value = 146
$ws1.Range("F37").Value = 4854
$ws1.Range("F38").Value = 21

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 627
$ws4.Range("F4").Value  = 1302
$ws4.Range("F5").Value  = 1170
$ws4.Range("F6").Value  = 14362
$ws4.Range("F7").Value  = 16658
$ws4.Range("F8").Value  = 15
$ws4.Range("F9").Value  = 107
$ws4.Range("F10").Value = 19
$ws4.Range("F11").Value = 50
$ws4.Range("F13").Value = 26
$ws4.Range("F14").Value = 51
$ws4.Range("F19").Value = 108
$ws4.Range("F20").Value = 38
$ws4.Range("F21").Value = 1271
$ws4.Range("F22").Value = 137
$ws4.Range("F25").Value = 42
$ws4.Range("F26").Value = 22
$ws4.Range("F28").Value = 6781
$ws4.Range("F30").Value = 23
$ws4.Range("F31").Value = 1125
$ws4.Range("F36").Value = 5766
$ws4.Range("F37").Value = 112
$ws4.Range("F38").Value = 146
$ws4.Range("F40").Value = 4854
$ws4.Range("F41").Value = 21
